# Bug fix in Eduati data files:
#  - Sheet1 ("HT115_noCTRL_meas") had leftover filler rows (45-87) that only
#    contained a running index in column A with no real measurement data.
#    Remove them so the sheet's used range matches the real data (A1:N44),
#    same as Sheet2/Sheet3.
#  - Make Sheet1 the active tab/selection (instead of Sheet3), and leave the
#    cursor positioned further down the sheet (D51) with the view scrolled
#    so row 36 is at the top.

$wb = $excel.ActiveWorkbook
$sheets = $wb.Worksheets

$ws1 = $sheets.Item(1)

# Drop the stray trailing rows (45:87) on Sheet1 that only held a bare index
# in column A - not part of the actual dataset.
$ws1.Rows("45:87").Delete() | Out-Null

# Make Sheet1 the active sheet/tab (previously Sheet3 was tabSelected and the
# workbook remembered Sheet3 as the active tab).
$ws1.Activate() | Out-Null

# Scroll the view down and leave the selection at D51, matching where the
# author was last working in the trimmed sheet.
$win = $excel.ActiveWindow
$win.ScrollRow = 36
$win.ScrollColumn = 1
$ws1.Range("D51").Select() | Out-Null

Write-Output "Trimmed Sheet1 to A1:N44 and made it the active sheet"
